$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Append the new log row (row 9)
$ws.Range("A9").Value = "Interne taak"
$ws.Range("B9").Value = "kwaliteit@testbedrijf123.nl"
$ws.Range("C9").Value = "Leg dit even neer bij Koen."
$ws.Range("D9").Value = "Intern verzoek / Actie voor medewerker"
$ws.Range("E9").Value = "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl."
$ws.Range("F9").Value = "2025-08-14 20:32:52"
$ws.Range("G9").Value = "Nee"
$ws.Range("H9").Value = "Ja"
$ws.Range("I9").Value = "Nee"
$ws.Range("J9").Value = "Nee"

# Extend the conditional-formatting ranges so they keep covering the new row
$ws.Range("D2:D8").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D9"))
$ws.Range("G2:G8").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G9"))
$ws.Range("H2:H8").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H9"))
$ws.Range("I2:I8").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I9"))
$ws.Range("J2:J8").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J9"))

# Update the Dashboard summary count for this category
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 8
